# Brain: Completing the Acceptance Criteria for Software Design Capability
#
# Fills in the (previously empty) "Acceptance Criteria" column (E) for the
# Software Design capability rows (121-136) and resizes those rows so the
# wrapped text displays fully, matching the author's manual edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> Acceptance Criteria text (column E)
$acceptanceCriteria = @{
    121 = 'Design and technical decisions are done in an adhoc manner with no vision of the application road map. Technical team members are not aware of the overall technical design decisions and this is not documented in any form.'
    122 = 'At this level of proficiency, the team does a lot of Big Up-Front designs for the systems as opposed to a just enough design approach. There are several ceremonial design sessions for the entire application without necessarily taken into effect the fact that there might need to be changes should the product vision change.'
    123 = 'Technologists at this level of proficiency design highly tightly coupled and brittle systems. The goal should be to start designing loosely coupled systems so that they are independently deployable and can be monitored in isolation.'
    124 = 'At this proficiency level, the designs are done more regularly and reviews are put in place to ensure that product vision change causes a validation of the technical design. The team also gets a view of the technical design decisions and roadmap.'
    125 = 'Design assumptions are tracked and recorded on the team wall and these are validated during the iteration as well as during the iteration planning sessions.'
    126 = 'Design issues identified are tracked and prioritised in the backlog and attended to similar to the technical debts.'
    127 = 'The system non-functional requirements are documented and tracked. These are also prioritised in the backlog.'
    128 = 'The team emphasizes a modular design of the application with clearly defined interfaces between the modules and clear communication protocols between the modules (for example clearly defined APIs using RESTful protocols).'
    129 = 'The technical design is owned by the feature teams. The technical architect is part of the core team and is not a stakeholder external to the team making decisions on their behalf.'
    130 = 'Some design decisions can be linked to the requirements, however, not all the design decisions are.'
    131 = 'The design requirements are captured in automated test cases, to ensure that these have been met.'
    132 = 'Non-functional requirements are clearly defined and documented. Associated metrics are tracked on the team dashboards automatically and are constantly monitored and optimized.'
    133 = 'Modules of the system can be simulated using stubs or mocks where applicable to ensure autonomy and isolation.'
    134 = 'Design is captured in a model and generated by software applications.'
    135 = 'All the design decisions can be traced to system requirements.'
    136 = 'Design happens as a Just-In-Time activity following paradigms like TDD and Evolutionary Architecture.'
}

# Row -> final row height (only rows whose wrapped text now needs more room)
$rowHeights = @{
    121 = 57.6
    122 = 72
    123 = 57.6
    124 = 57.6
    125 = 43.2
    126 = 28.8
    127 = 28.8
    128 = 57.6
    129 = 43.2
    130 = 28.8
    131 = 28.8
    132 = 43.2
    133 = 28.8
    136 = 28.8
}

foreach ($row in 121..136) {
    $ws.Range("E$row").Value = $acceptanceCriteria[$row]
    if ($rowHeights.ContainsKey($row)) {
        $ws.Rows($row).RowHeight = $rowHeights[$row]
    }
}

# Reflect where the author ended up looking (scrolled / selected) after
# finishing the edits.
$ws.Range("E137").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 126
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
